$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.174.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.052.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.664"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.95%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.21"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.29%  "
$ws.Range("E9").Value = "  -2.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0782"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.36%  "
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.13"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.885"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.351.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.73"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.043.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +17.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.161.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0895"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.60%  "
$ws.Range("E21").Value = "  -1.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "237.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.91%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  +3.26%  "
$ws.Range("E25").Value = "  +2.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "169.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.32%  "
$ws.Range("E29").Value = "  -1.17%  "
$ws.Range("E30").Value = "  -0.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0621"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.50"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0890"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.67%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  -2.84%  "
$ws.Range("E37").Value = "  +0.58%  "
$ws.Range("E38").Value = "  -2.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.30"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +14.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.12"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1000"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -14.48%  "
$ws.Range("E42").Value = "  -2.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "96.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.17%  "
$ws.Range("E46").Value = "  -2.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.269.66"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.80%  "
$ws.Range("E48").Value = "  -3.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.83"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.239.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.68%  "
